# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13343
$ws1.Range("F4").Value = 647
$ws1.Range("F5").Value = 216
$ws1.Range("F6").Value = 446
$ws1.Range("F7").Value = 1317
$ws1.Range("F8").Value = 125

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13343
$ws4.Range("F4").Value = 647
$ws4.Range("F5").Value = 216
$ws4.Range("F8").Value = 446
$ws4.Range("F9").Value = 1317
$ws4.Range("F11").Value = 125
